$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.197.23"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.833.08"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.39"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6650"
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07346"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2913"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.57"
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07683"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.842.67"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.965"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6648"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.64"
$ws.Range("E15").Value = "  -4.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.074"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.179.25"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008266"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.42"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "224.91"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.119"
$ws.Range("E22").Value = "  -3.98%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.46"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.619"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1391"
$ws.Range("E26").Value = "  -4.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.88"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.511"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.104"
$ws.Range("E29").Value = "  -4.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.027"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.191"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05279"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.863"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7458"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.685"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.312.32"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01793"
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9148"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.08556"
$ws.Range("E41").Value = "  +19.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.939"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.00000000132"
$ws.Range("E44").Value = "  +7.50%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.79"
$ws.Range("E45").Value = "  -3.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.970.41"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5163"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.759"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.13"
$ws.Range("E49").Value = "  -3.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.047"
$ws.Range("E50").Value = "  -5.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05929"
$ws.Range("E51").Value = "  -0.48%  "
